$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style/format of column A (row 11) down to the newly added rows (12-38)
$ws.Range("A11").Copy() | Out-Null
$ws.Range("A12:A38").PasteSpecial(-4122) | Out-Null

# Update / populate cell values for rows 2-38
$ws.Range("A2").Value = 44
$ws.Range("B2").Value = "U931101109048"
$ws.Range("A3").Value = 42
$ws.Range("B3").Value = "M931235210024"
$ws.Range("A4").Value = 39
$ws.Range("B4").Value = "Z931100609006"
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = "L394201008038"
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "G931101109060"
$ws.Range("A7").Value = 38
$ws.Range("B7").Value = "W931100608061"
$ws.Range("A8").Value = 27
$ws.Range("B8").Value = "F931235210018"
$ws.Range("A9").Value = 43
$ws.Range("B9").Value = "A931252110030"
$ws.Range("A10").Value = 6
$ws.Range("B10").Value = "M931100609016"
$ws.Range("A11").Value = 16
$ws.Range("B11").Value = "Q931235212001"
$ws.Range("A12").Value = 40
$ws.Range("B12").Value = "N931325209054"
$ws.Range("A13").Value = 32
$ws.Range("B13").Value = "Z931383214002"
$ws.Range("A14").Value = 19
$ws.Range("B14").Value = "F931100509027"
$ws.Range("A15").Value = 23
$ws.Range("B15").Value = "R931253116053"
$ws.Range("A16").Value = 12
$ws.Range("B16").Value = "K931100609063"
$ws.Range("A17").Value = 11
$ws.Range("B17").Value = "V931414517045"
$ws.Range("A18").Value = 45
$ws.Range("B18").Value = "V931101109041"
$ws.Range("A19").Value = 21
$ws.Range("B19").Value = "U931101109019"
$ws.Range("A20").Value = 4
$ws.Range("B20").Value = "C931252508049"
$ws.Range("A21").Value = 5
$ws.Range("B21").Value = "V802229210007"
$ws.Range("A22").Value = 7
$ws.Range("B22").Value = "N931100609007"
$ws.Range("A23").Value = 9
$ws.Range("B23").Value = "E931252916073"
$ws.Range("A24").Value = 15
$ws.Range("B24").Value = "M931252710007"
$ws.Range("A25").Value = 36
$ws.Range("B25").Value = "J931384210007"
$ws.Range("A26").Value = 35
$ws.Range("B26").Value = "V931101109012"
$ws.Range("A27").Value = 34
$ws.Range("B27").Value = "K931101109004"
$ws.Range("A28").Value = 33
$ws.Range("B28").Value = "C931100609037"
$ws.Range("A29").Value = 17
$ws.Range("B29").Value = "Y931252211003"
$ws.Range("A30").Value = 30
$ws.Range("B30").Value = "N931383610018"
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = "G931383410017"
$ws.Range("A32").Value = 28
$ws.Range("B32").Value = "T931252911047"
$ws.Range("A33").Value = 13
$ws.Range("B33").Value = "M931252909052"
$ws.Range("A34").Value = 26
$ws.Range("B34").Value = "F931100609041"
$ws.Range("A35").Value = 25
$ws.Range("B35").Value = "X931252710015"
$ws.Range("A36").Value = 14
$ws.Range("B36").Value = "W931321110033"
$ws.Range("A37").Value = 20
$ws.Range("B37").Value = "C931321610014"
$ws.Range("A38").Value = 18
$ws.Range("B38").Value = "J931101109071"
